$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '60.784.49'
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").Value = '3.383.70'
$ws.Range("E3").Value = '  -1.95%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.95%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '3.380.54'
$ws.Range("E8").Value = '  -2.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.473'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.124'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("E12").Value = '  +1.81%  '
$ws.Range("D13").Value = '3.959.62'
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("E16").Value = '  -2.26%  '
$ws.Range("D17").Value = '3.386.10'
$ws.Range("E17").Value = '  -1.87%  '
$ws.Range("D18").Value = '60.857.07'
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("E20").Value = '  +0.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '387.61'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.562'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.95'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.996'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000120'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.07%  '
$ws.Range("D27").Value = '3.538.62'
$ws.Range("E27").Value = '  -1.31%  '
$ws.Range("E28").Value = '  +0.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.33'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.56%  '
$ws.Range("E32").Value = '  -7.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.16'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.36%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.78'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.99'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("D37").Value = '3.410.38'
$ws.Range("E37").Value = '  -1.74%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '167.41'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.96%  '
$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.06'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.58%  '
$ws.Range("E40").Value = '  -1.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0781'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.88'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.788'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.89%  '
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("E46").Value = '  -0.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.69'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.14%  '
$ws.Range("D48").Value = '2.552.61'
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.12'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.89'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.97'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.49%  '
